$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit performs a cyclic rotation of the species-observation data among
# rows 2, 3 and 4 (row2 <- old row4, row3 <- old row2, row4 <- old row3),
# while shared/common columns (C, D, I, P, S, T, U, V, W, Y, Z, AA, AB, AD,
# AE, AG, AT, AW, AX, AY) remain identical across the three rows and are
# therefore unaffected.

# Row 2 (new values, previously held by row 4)
$ws.Range("A2").Value = 62689594
$ws.Range("B2").Value = 90319
$ws.Range("E2").Value = 4769
$ws.Range("F2").Value = "Svavelriska"
$ws.Range("G2").Value = "Lactarius scrobiculatus"
$ws.Range("H2").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q2").Value = 496603.0244192505
$ws.Range("R2").Value = 6593644.236533851
$ws.Range("AI2").Value = "Lövrik barrskog"
$ws.Range("AO2").Value = "Fuktig mark"

# Row 3 (new values, previously held by row 2)
$ws.Range("A3").Value = 62689582
$ws.Range("B3").Value = 93054
$ws.Range("E3").Value = 2810
$ws.Range("F3").Value = "Västlig hakmossa"
$ws.Range("G3").Value = "Rhytidiadelphus loreus"
$ws.Range("H3").Value = "(Hedw.) Warnst."
$ws.Range("Q3").Value = 496491.2429489095
$ws.Range("R3").Value = 6593829.754447529
$ws.Range("AI3").Value = "Granskog"
$ws.Range("AO3").Value = "Marken"

# Row 4 (new values, previously held by row 3)
$ws.Range("A4").Value = 62689554
$ws.Range("B4").Value = 93044
$ws.Range("E4").Value = 2809
$ws.Range("F4").Value = "Mörk husmossa"
$ws.Range("G4").Value = "Hylocomiastrum umbratum"
$ws.Range("H4").Value = "(Hedw.) M.Fleisch."
$ws.Range("Q4").Value = 496589.8900309857
$ws.Range("R4").Value = 6593750.926901286
$ws.Range("AI4").Value = "Lövrik barrskog"
$ws.Range("AO4").Value = "Marken"
